$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data: each row 31-38 now holds what used to be in row 32-39,
# and a brand-new week is appended as row 39.

$data = @(
    @{ Row = 31; D = 45077; I = "Primera"; J = 60;  K = 8000;  L = 8000;  M = 8000;  N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 133; Q = 60 },
    @{ Row = 32; D = 45077; I = "Segunda"; J = 60;  K = 6000;  L = 6000;  M = 6000;  N = "$/caja 90 unidades"; O = "Región de Arica y Parinacota"; P = 67;  Q = 90 },
    @{ Row = 33; D = 44979; I = "Primera"; J = 50;  K = 9000;  L = 9000;  M = 9000;  N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 150; Q = 60 },
    @{ Row = 34; D = 45063; I = "Primera"; J = 50;  K = 8500;  L = 8500;  M = 8500;  N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 142; Q = 60 },
    @{ Row = 35; D = 44162; I = "Primera"; J = 43;  K = 8000;  L = 8500;  M = 8209;  N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 137; Q = 60 },
    @{ Row = 36; D = 44776; I = "Primera"; J = 60;  K = 11000; L = 12000; M = 11500; N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 192; Q = 60 },
    @{ Row = 37; D = 44671; I = "Primera"; J = 160; K = 6000;  L = 7000;  M = 6500;  N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 108; Q = 60 },
    @{ Row = 38; D = 44202; I = "Primera"; J = 50;  K = 8000;  L = 9000;  M = 8400;  N = "$/caja 60 unidades"; O = "Región del Maule";             P = 140; Q = 60 },
    @{ Row = 39; D = 45117; I = "Primera"; J = 40;  K = 10000; L = 10000; M = 10000; N = "$/caja 60 unidades"; O = "Región de Arica y Parinacota"; P = 167; Q = 60 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("D$r").Value = $item.D
    $ws.Range("I$r").Value = $item.I
    $ws.Range("J$r").Value = $item.J
    $ws.Range("K$r").Value = $item.K
    $ws.Range("L$r").Value = $item.L
    $ws.Range("M$r").Value = $item.M
    $ws.Range("N$r").Value = $item.N
    $ws.Range("O$r").Value = $item.O
    $ws.Range("P$r").Value = $item.P
    $ws.Range("Q$r").Value = $item.Q
}
